$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 82; this shifts the existing rows 82-98 down to 83-99
# and keeps their data/formatting intact.
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with the new weekly price record.
$ws.Cells.Item(82, 1).Value = 2
$ws.Cells.Item(82, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(82, 3).Value = "Coquimbo"
$ws.Cells.Item(82, 4).Value = 45127
$ws.Cells.Item(82, 5).Value = 4
$ws.Cells.Item(82, 6).Value = 100112026
$ws.Cells.Item(82, 7).Value = "Haba"
$ws.Cells.Item(82, 8).Value = "Sin especificar"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 700
$ws.Cells.Item(82, 11).Value = 10000
$ws.Cells.Item(82, 12).Value = 11000
$ws.Cells.Item(82, 13).Value = 10500
$ws.Cells.Item(82, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(82, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(82, 16).Value = 420
$ws.Cells.Item(82, 17).Value = 25
$ws.Cells.Item(82, 18).Value = "Hortaliza"
